$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-GuidRow($row, $guid, $info) {
    $cellA = $ws.Cells.Item($row, 1)
    # Force the GUID code (e.g. "000083") to be stored as text rather than
    # being auto-converted to a number, matching the existing rows above.
    $cellA.NumberFormat = "@"
    $cellA.Value = $guid

    # Re-apply the plain (General/default) style used by the rows above,
    # so the new cell matches their formatting exactly instead of keeping
    # a leftover "@" text-format style.
    $ws.Cells.Item(83, 1).Copy() | Out-Null
    $cellA.PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 2).Value = $info
}

Set-GuidRow 84 "000083" "Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 08-Mar-2023 10:28:53"

Set-GuidRow 85 "000084" "Details: IMU CJMCU-20948 Raw Data Reading - Raw Plots of 100 points of Accel Gyro and Magnometer data separately. Script used: Read_IMU.  Dataset used: Arduino Serial Output of IMU CJMCU-20948. File Location: Visualisations/IMU_RealRawData. Date Generated: 08-Mar-2023 10:29:27"

Set-GuidRow 86 "000085" "Details: Baton Tip Pose Transformation. IMU CJMCU-20948 Data Reading, Fused with imufilter, transformed with BatonTip_Transformation. Script used: BatonTipPoseVisualisation.  Dataset used: IMU data: IMU_Orientation_Reading_08_03_23. Transformed Baton tip data: BatonTipPose_08_03_23.. File Location: Visualisations/IMU_TransformedBatonTipPose. Date Generated: 08-Mar-2023 11:41:21"

Set-GuidRow 87 "000086" "Details: Baton Tip Pose Transformation. IMU CJMCU-20948 Data Reading, Fused with imufilter, transformed with BatonTip_Transformation. Script used: BatonTipPoseVisualisation.  Dataset used: IMU data: IMU_Orientation_Reading_08_03_23. Transformed Baton tip data: BatonTipPose_08_03_23.. File Location: Visualisations/IMU_TransformedBatonTipPose. Date Generated: 08-Mar-2023 11:44:01"
